$d = $word.ActiveDocument

# Order matters: "884÷4=" becomes "523÷2=", which collides with the
# pre-existing "523÷2=" cell (which itself becomes "567÷7="). Perform the
# "523÷2=" -> "567÷7=" replacement first so the later "884÷4=" -> "523÷2="
# replacement does not get clobbered by its own output being re-matched.

$d.Content.Find.Execute("2025-07-12 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-13 Sunday", 2) | Out-Null
$d.Content.Find.Execute("832÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "153÷8=", 2) | Out-Null
$d.Content.Find.Execute("469÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "385÷7=", 2) | Out-Null
$d.Content.Find.Execute("940÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "253÷7=", 2) | Out-Null
$d.Content.Find.Execute("977÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "600÷9=", 2) | Out-Null
$d.Content.Find.Execute("178÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "999÷5=", 2) | Out-Null
$d.Content.Find.Execute("100÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "168÷4=", 2) | Out-Null
$d.Content.Find.Execute("827÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "919÷6=", 2) | Out-Null
$d.Content.Find.Execute("513÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "873÷2=", 2) | Out-Null
$d.Content.Find.Execute("646÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "986÷2=", 2) | Out-Null
$d.Content.Find.Execute("923÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "199÷6=", 2) | Out-Null
$d.Content.Find.Execute("186÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "565÷2=", 2) | Out-Null
$d.Content.Find.Execute("990÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "947÷5=", 2) | Out-Null
$d.Content.Find.Execute("530÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "794÷6=", 2) | Out-Null
$d.Content.Find.Execute("986÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "781÷3=", 2) | Out-Null
$d.Content.Find.Execute("926÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "148÷7=", 2) | Out-Null
$d.Content.Find.Execute("526÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "205÷5=", 2) | Out-Null
$d.Content.Find.Execute("523÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "567÷7=", 2) | Out-Null
$d.Content.Find.Execute("844÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "259÷9=", 2) | Out-Null
$d.Content.Find.Execute("223÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "828÷9=", 2) | Out-Null
$d.Content.Find.Execute("133÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "691÷5=", 2) | Out-Null
$d.Content.Find.Execute("382÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "268÷2=", 2) | Out-Null
$d.Content.Find.Execute("135÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "829÷5=", 2) | Out-Null
$d.Content.Find.Execute("608÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "713÷7=", 2) | Out-Null
$d.Content.Find.Execute("995÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "297÷3=", 2) | Out-Null
$d.Content.Find.Execute("884÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "523÷2=", 2) | Out-Null

$d.Save()
